# "Fix Rappel AV cloture issue" — the underlying "état des taxes" report was
# regenerated: the contract grouped under "905/LF/TADLA OUARDIGHA ZAYANE"
# (rows 2-7) now splits its "Rappel" (MT brut/Taxe) amounts differently
# between the "loyer" columns (H/J) and the "Rappel" columns (L/M), several
# contracts further down were re-keyed/reordered against updated CIN/IF and
# montant data, and the trailing 3 rows (two "Direction régionale" rows plus
# the totals row) were dropped — the totals row that used to live at row 20
# now lives at row 17 with recomputed totals. Net effect: sheet dimension
# shrinks from A1:O20 to A1:O17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three trailing rows (old rows 18, 19, 20); rows 2-17 keep their
# addresses since everything removed sits below them.
$ws.Range("A18:O20").EntireRow.Delete()

# --- Rows 2-4: "905/LF/TADLA OUARDIGHA ZAYANE" / Logement de fonction -----
# Rappel split moves from L/M into H/J ("--" placeholders become real loyer
# amounts), and the Rappel (L/M) + net (O) amounts drop accordingly.
$ws.Range("H2").Value = 3333.33
$ws.Range("J2").Value = 333.33
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("O2").Value = 3000

$ws.Range("H3").Value = 3333.33
$ws.Range("J3").Value = 333.33
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 3000

$ws.Range("H4").Value = 3333.33
$ws.Range("J4").Value = 333.33
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("O4").Value = 3000

# --- Row 5: now "052/FKIH BEN SALEH/AV1" / ZERNAKH ABDELLAH (was a third
# "905/LF.../NASIRI HASNAA" row) ---------------------------------------
$ws.Range("A5").Value = "052/FKIH BEN SALEH/AV1"
$ws.Range("B5").Value = "Point de vente"
$ws.Range("C5").Value = "IB19558"
$ws.Range("D5").Value = "ZERNAKH ABDELLAH"
$ws.Range("E5").Value = "oui"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "--"
$ws.Range("J5").Value = "--"
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 450
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 2550

# --- Row 6: second "052/FKIH BEN SALEH/AV1" / ZERNAKH ABDELLAH row -------
$ws.Range("A6").Value = "052/FKIH BEN SALEH/AV1"
$ws.Range("B6").Value = "Point de vente"
$ws.Range("C6").Value = "IB19558"
$ws.Range("D6").Value = "ZERNAKH ABDELLAH"
$ws.Range("E6").Value = "oui"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 12000
$ws.Range("J6").Value = 0
$ws.Range("N6").Value = "--"
$ws.Range("O6").Value = 12000

# --- Row 7: now "605/KHOURIBGA NAHDA" / MOHAMED BADRANE ------------------
$ws.Range("A7").Value = "605/KHOURIBGA NAHDA"
$ws.Range("B7").Value = "Point de vente"
$ws.Range("C7").Value = "I83603"
$ws.Range("D7").Value = "MOHAMED BADRANE"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 446.42
$ws.Range("J7").Value = 0
$ws.Range("N7").Value = "--"
$ws.Range("O7").Value = 446.42

# --- Row 8: JEMAA HORMI, montant changes to 1250 --------------------------
$ws.Range("C8").Value = "B219321"
$ws.Range("D8").Value = "JEMAA HORMI"
$ws.Range("H8").Value = 1250
$ws.Range("O8").Value = 1250

# --- Row 9: DOUNIA LAMKADDAM, montant changes to 937.5 --------------------
$ws.Range("C9").Value = "BK646476"
$ws.Range("D9").Value = "DOUNIA LAMKADDAM"
$ws.Range("H9").Value = 937.5
$ws.Range("O9").Value = 937.5

# --- Row 10: SOFIA BADRANE, now taxed at 10% / montant 3750 ---------------
$ws.Range("C10").Value = "CIN605"
$ws.Range("D10").Value = "SOFIA BADRANE"
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = 3750
$ws.Range("J10").Value = 375
$ws.Range("O10").Value = 3375

# --- Row 11: LATIFA BADRANE, montant 223.21 -------------------------------
$ws.Range("C11").Value = "I150156"
$ws.Range("D11").Value = "LATIFA BADRANE"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 223.21
$ws.Range("J11").Value = 0
$ws.Range("O11").Value = 223.21

# --- Row 12: NADIA BADRANE (montant unchanged, only CIN/name shift) ------
$ws.Range("C12").Value = "B171710"
$ws.Range("D12").Value = "NADIA BADRANE"

# --- Row 13: OUAFA BADRANE (montant unchanged, only CIN/name shift) ------
$ws.Range("C13").Value = "Q194939"
$ws.Range("D13").Value = "OUAFA BADRANE"

# --- Row 14: SAID BADRANE, montant changes to 446.45 ----------------------
$ws.Range("C14").Value = "I210578"
$ws.Range("D14").Value = "SAID BADRANE"
$ws.Range("H14").Value = 446.45
$ws.Range("O14").Value = 446.45

# --- Row 15: now "905/TADLA OUARDIGHA ZAYANE" / Direction régionale /
# NOUBAIL MOUNTASSIR (was "605/KHOURIBGA NAHDA" / SAID BADRANE) ----------
$ws.Range("A15").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("B15").Value = "Direction régionale"
$ws.Range("C15").Value = "Q251990"
$ws.Range("D15").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 6750
$ws.Range("J15").Value = 675
$ws.Range("O15").Value = 6075

# --- Row 16: now "905/TADLA OUARDIGHA ZAYANE" / Direction régionale /
# NOUBAIL MOHAMMED (was "052/FKIH BEN SALEH" / ZERNAKH ABDELLAH) ---------
$ws.Range("A16").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("B16").Value = "Direction régionale"
$ws.Range("C16").Value = "IR801997"
$ws.Range("D16").Value = "NOUBAIL MOHAMMED"
$ws.Range("E16").Value = "non"
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 6750
$ws.Range("J16").Value = 675
$ws.Range("O16").Value = 6075

# --- Row 17: the report's trailing totals row (used to be row 20); the
# label cells stay blank ("A17:G17" = " "), only the totals changed. -----
$ws.Range("A17").Value = " "
$ws.Range("B17").Value = " "
$ws.Range("C17").Value = " "
$ws.Range("D17").Value = " "
$ws.Range("E17").Value = " "
$ws.Range("F17").Value = " "
$ws.Range("G17").Value = " "
$ws.Range("H17").Value = 42999.99
$ws.Range("J17").Value = 2724.99
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 450
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 42825
